# Generate Report for Handback
# Update the handback/handoff timestamps and status for the second file
# (b6ffff3b-a984-4897-944a-f55dee0f69ba) now that a new handback round has
# completed, and roll the "Latest HO Xliff Generate Date" on the Overview
# sheet up to match.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-24 20:50:21"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-24 20:50:00"
$wsZhCn.Range("K3").Value = "2016-08-24 20:50:38"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-24 20:50:21"
$wsDeDe.Range("K3").Value = "2016-08-24 20:50:46"
